$d = $word.ActiveDocument

# Update the date/weekday heading
$d.Content.Find.Execute("2026-02-07 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-08 Sunday", 2)

# Update the division problems in the table. The table has 20 rows x 5 columns,
# but only every 4th row (1, 5, 9, 13, 17) actually contains text.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "37÷8="
$t.Cell(1, 2).Range.Text  = "31÷2="
$t.Cell(1, 3).Range.Text  = "86÷5="
$t.Cell(1, 4).Range.Text  = "29÷9="
$t.Cell(1, 5).Range.Text  = "22÷4="

$t.Cell(5, 1).Range.Text  = "85÷9="
$t.Cell(5, 2).Range.Text  = "58÷9="
$t.Cell(5, 3).Range.Text  = "14÷9="
$t.Cell(5, 4).Range.Text  = "74÷7="
$t.Cell(5, 5).Range.Text  = "96÷9="

$t.Cell(9, 1).Range.Text  = "26÷3="
$t.Cell(9, 2).Range.Text  = "56÷6="
$t.Cell(9, 3).Range.Text  = "64÷9="
$t.Cell(9, 4).Range.Text  = "21÷5="
$t.Cell(9, 5).Range.Text  = "59÷5="

$t.Cell(13, 1).Range.Text = "78÷3="
$t.Cell(13, 2).Range.Text = "63÷8="
$t.Cell(13, 3).Range.Text = "30÷6="
$t.Cell(13, 4).Range.Text = "68÷7="
$t.Cell(13, 5).Range.Text = "80÷3="

$t.Cell(17, 1).Range.Text = "29÷4="
$t.Cell(17, 2).Range.Text = "78÷7="
$t.Cell(17, 3).Range.Text = "64÷4="
$t.Cell(17, 4).Range.Text = "27÷6="
$t.Cell(17, 5).Range.Text = "84÷7="
